$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.774.90"
$ws.Range("E2").Value = "  +4.00%  "

# Row 3
$ws.Range("D3").Value = "3.086.65"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'236.47"
$ws.Range("E5").Value = "  +1.37%  "

# Row 6
$ws.Range("D6").Value = "'603.41"

# Row 7
$ws.Range("E7").Value = "  +2.77%  "

# Row 8
$ws.Range("D8").Value = "'0.378"
$ws.Range("E8").Value = "  -1.35%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").Value = "3.081.60"
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$ws.Range("D11").Value = "'0.780"
$ws.Range("E11").Value = "  +1.73%  "

# Row 12
$ws.Range("E12").Value = "  -0.40%  "

# Row 13
$ws.Range("D13").Value = "95.007.85"
$ws.Range("E13").Value = "  +3.39%  "

# Row 14
$ws.Range("E14").Value = "  -2.45%  "

# Row 15
$ws.Range("E15").Value = "  -0.89%  "

# Row 16
$ws.Range("D16").Value = "'5.28"
$ws.Range("E16").Value = "  -1.78%  "

# Row 17
$ws.Range("D17").Value = "3.650.64"
$ws.Range("E17").Value = "  -0.76%  "

# Row 18
$ws.Range("D18").Value = "3.056.48"
$ws.Range("E18").Value = "  -1.02%  "

# Row 19
$ws.Range("D19").Value = "'3.49"
$ws.Range("E19").Value = "  -8.16%  "

# Row 20
$ws.Range("D20").Value = "'14.18"
$ws.Range("E20").Value = "  -1.45%  "

# Row 21
$ws.Range("D21").Value = "'449.11"
$ws.Range("E21").Value = "  +3.35%  "

# Row 22
$ws.Range("D22").Value = "'5.57"
$ws.Range("E22").Value = "  -3.14%  "

# Row 23
$ws.Range("E23").Value = "  -3.70%  "

# Row 24
$ws.Range("D24").Value = "'8.61"
$ws.Range("E24").Value = "  -4.64%  "

# Row 25
$ws.Range("D25").Value = "'5.45"
$ws.Range("E25").Value = "  -2.00%  "

# Row 26
$ws.Range("D26").Value = "'84.62"
$ws.Range("E26").Value = "  -0.41%  "

# Row 27
$ws.Range("D27").Value = "'11.47"
$ws.Range("E27").Value = "  +1.92%  "

# Row 28
$ws.Range("D28").Value = "3.233.77"
$ws.Range("E28").Value = "  -0.77%  "

# Row 29
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.240"
$ws.Range("E30").Value = "  +4.34%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.128"
$ws.Range("E31").Value = "  +2.21%  "

# Row 32
$ws.Range("E32").Value = "  +0.32%  "

# Row 33
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").Value = "'8.83"
$ws.Range("E34").Value = "  -2.74%  "

# Row 35
$ws.Range("D35").Value = "'25.37"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").Value = "'7.24"
$ws.Range("E36").Value = "  -7.79%  "

# Row 37
$ws.Range("D37").Value = "'0.148"
$ws.Range("E37").Value = "  -4.51%  "

# Row 38
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'24.12"
$ws.Range("E38").Value = "  +1.18%  "

# Row 39
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'479.18"
$ws.Range("E39").Value = "  +2.80%  "

# Row 40
$ws.Range("B40").Value = "PancakeSwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D40").Value = "'1.85"
$ws.Range("E40").Value = "  -1.34%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.428"
$ws.Range("E41").Value = "  -0.42%  "

# Row 42
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").Value = "'3.65"
$ws.Range("E42").Value = "  -6.06%  "

# Row 43
$ws.Range("E43").Value = "  -4.27%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'3.09"
$ws.Range("E45").Value = "  -5.05%  "

# Row 46
$ws.Range("D46").Value = "'160.87"
$ws.Range("E46").Value = "  +0.43%  "

# Row 47
$ws.Range("D47").Value = "'0.671"
$ws.Range("E47").Value = "  -0.87%  "

# Row 48
$ws.Range("E48").Value = "  -0.97%  "

# Row 49
$ws.Range("D49").Value = "'0.000275"
$ws.Range("E49").Value = "  +13.81%  "

# Row 50
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.08%  "
